$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.137.81"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "3.028.05"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.25"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.33"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.024.47"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -4.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.50"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "3.514.17"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").Value = "62.182.07"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "3.025.84"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.70"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.62"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.76"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.21"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.86"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.33"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.35"
$ws.Range("E35").Value = "  -4.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "461.00"
$ws.Range("E37").Value = "  -8.30%  "
$ws.Range("D38").Value = "3.237.87"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0801"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0387"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.18"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.49"
$ws.Range("E43").Value = "  -7.55%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.95"
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.246"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.109"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.84"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "0.0₃0500"
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("E51").Value = "  +7.20%  "
